# ==========================================================================
# Adds a new "2022-Q4" sheet (fund holdings snapshot) right after "总计",
# and inserts a corresponding summary row into "总计" itself.
# ==========================================================================

function Set-TextValue($sheet, $row, $col, $val) {
    # Excel auto-detects numeric-looking strings and stores them as numbers,
    # which would strip things like leading zeros in fund codes or trailing
    # zeros in percentages. A leading apostrophe forces genuine text entry,
    # exactly like a user typing '010723 into a cell.
    $sheet.Cells.Item($row, $col).Value = "'" + $val
}

$wb = $excel.ActiveWorkbook
$zj = $wb.Worksheets.Item(1)          # "总计" summary sheet (always sheet #1)
$q3 = $wb.Worksheets.Item("2022-Q3")  # template to clone for the new quarter

# ---- 1. Create the "2022-Q4" sheet by duplicating "2022-Q3" -------------
# Copying preserves every style (header bold+border, index-column style,
# column widths, etc.) exactly, and placing it "after" 总计 puts it in the
# correct tab position in one step.
$q3.Copy($null, $zj)
$q4 = $wb.Worksheets.Item(2)
$q4.Name = "2022-Q4"

# 2022-Q3 has 42 data rows (rows 2-43); 2022-Q4 only has 40 (rows 2-41), so
# trim the two trailing rows inherited from the clone.
$q4.Rows.Item(42).Delete()
$q4.Rows.Item(42).Delete()

# ---- 2. Overwrite the cloned sheet's data with the 2022-Q4 figures ------
# (Column A - the 0-based row index - already reads 0..39 after the clone +
# trim above, so it does not need to be touched.)

Set-TextValue $q4 2 2 "010723"
$q4.Cells.Item(2, 3).Value = "中欧价值成长混合A"
Set-TextValue $q4 2 4 "21.90"
Set-TextValue $q4 2 5 "83.01"
Set-TextValue $q4 2 6 "4.48"
Set-TextValue $q4 2 7 "0.9811"
$q4.Cells.Item(2, 8).Value = 4

Set-TextValue $q4 3 2 "166009"
$q4.Cells.Item(3, 3).Value = "中欧新动力混合（LOF）A"
Set-TextValue $q4 3 4 "15.62"
Set-TextValue $q4 3 5 "86.77"
Set-TextValue $q4 3 6 "4.81"
Set-TextValue $q4 3 7 "0.7513"
$q4.Cells.Item(3, 8).Value = 3

Set-TextValue $q4 4 2 "009210"
$q4.Cells.Item(4, 3).Value = "中欧嘉和三年持有期混合A"
Set-TextValue $q4 4 4 "14.41"
Set-TextValue $q4 4 5 "93.03"
Set-TextValue $q4 4 6 "4.57"
Set-TextValue $q4 4 7 "0.6585"
$q4.Cells.Item(4, 8).Value = 2

Set-TextValue $q4 5 2 "519692"
$q4.Cells.Item(5, 3).Value = "交银成长混合A"
Set-TextValue $q4 5 4 "25.08"
Set-TextValue $q4 5 5 "80.90"
Set-TextValue $q4 5 6 "2.56"
Set-TextValue $q4 5 7 "0.6420"
$q4.Cells.Item(5, 8).Value = 10

Set-TextValue $q4 6 2 "010678"
$q4.Cells.Item(6, 3).Value = "中欧均衡成长混合A"
Set-TextValue $q4 6 4 "13.63"
Set-TextValue $q4 6 5 "83.93"
Set-TextValue $q4 6 6 "4.71"
Set-TextValue $q4 6 7 "0.6420"
$q4.Cells.Item(6, 8).Value = 4

Set-TextValue $q4 7 2 "519700"
$q4.Cells.Item(7, 3).Value = "交银主题优选混合A"
Set-TextValue $q4 7 4 "26.24"
Set-TextValue $q4 7 5 "73.61"
Set-TextValue $q4 7 6 "1.90"
Set-TextValue $q4 7 7 "0.4986"
$q4.Cells.Item(7, 8).Value = 9

Set-TextValue $q4 8 2 "519694"
$q4.Cells.Item(8, 3).Value = "交银蓝筹混合"
Set-TextValue $q4 8 4 "17.19"
Set-TextValue $q4 8 5 "81.14"
Set-TextValue $q4 8 6 "2.53"
Set-TextValue $q4 8 7 "0.4349"
$q4.Cells.Item(8, 8).Value = 10

Set-TextValue $q4 9 2 "010947"
$q4.Cells.Item(9, 3).Value = "中欧嘉选混合A"
Set-TextValue $q4 9 4 "12.93"
Set-TextValue $q4 9 5 "83.49"
Set-TextValue $q4 9 6 "3.22"
Set-TextValue $q4 9 7 "0.4163"
$q4.Cells.Item(9, 8).Value = 4

Set-TextValue $q4 10 2 "005421"
$q4.Cells.Item(10, 3).Value = "中欧嘉泽灵活配置混合"
Set-TextValue $q4 10 4 "6.80"
Set-TextValue $q4 10 5 "93.21"
Set-TextValue $q4 10 6 "5.34"
Set-TextValue $q4 10 7 "0.3631"
$q4.Cells.Item(10, 8).Value = 2

Set-TextValue $q4 11 2 "320001"
$q4.Cells.Item(11, 3).Value = "诺安平衡混合"
Set-TextValue $q4 11 4 "10.86"
Set-TextValue $q4 11 5 "71.09"
Set-TextValue $q4 11 6 "3.04"
Set-TextValue $q4 11 7 "0.3301"
$q4.Cells.Item(11, 8).Value = 10

Set-TextValue $q4 12 2 "009564"
$q4.Cells.Item(12, 3).Value = "汇安消费龙头混合A"
Set-TextValue $q4 12 4 "7.95"
Set-TextValue $q4 12 5 "94.58"
Set-TextValue $q4 12 6 "3.27"
Set-TextValue $q4 12 7 "0.2600"
$q4.Cells.Item(12, 8).Value = 9

Set-TextValue $q4 13 2 "013993"
$q4.Cells.Item(13, 3).Value = "中欧光熠一年持有期混合型证券投资基金A"
Set-TextValue $q4 13 4 "6.13"
Set-TextValue $q4 13 5 "87.20"
Set-TextValue $q4 13 6 "3.93"
Set-TextValue $q4 13 7 "0.2409"
$q4.Cells.Item(13, 8).Value = 3

Set-TextValue $q4 14 2 "013884"
$q4.Cells.Item(14, 3).Value = "交银主题优选混合C"
Set-TextValue $q4 14 4 "11.57"
Set-TextValue $q4 14 5 "73.61"
Set-TextValue $q4 14 6 "1.90"
Set-TextValue $q4 14 7 "0.2198"
$q4.Cells.Item(14, 8).Value = 9

Set-TextValue $q4 15 2 "011708"
$q4.Cells.Item(15, 3).Value = "中欧嘉益一年混合A"
Set-TextValue $q4 15 4 "4.37"
Set-TextValue $q4 15 5 "91.19"
Set-TextValue $q4 15 6 "4.23"
Set-TextValue $q4 15 7 "0.1849"
$q4.Cells.Item(15, 8).Value = 3

Set-TextValue $q4 16 2 "004236"
$q4.Cells.Item(16, 3).Value = "中欧新动力混合（LOF）C"
Set-TextValue $q4 16 4 "3.50"
Set-TextValue $q4 16 5 "86.77"
Set-TextValue $q4 16 6 "4.81"
Set-TextValue $q4 16 7 "0.1684"
$q4.Cells.Item(16, 8).Value = 3

Set-TextValue $q4 17 2 "012202"
$q4.Cells.Item(17, 3).Value = "中加消费优选混合A"
Set-TextValue $q4 17 4 "3.65"
Set-TextValue $q4 17 5 "88.36"
Set-TextValue $q4 17 6 "3.13"
Set-TextValue $q4 17 7 "0.1142"
$q4.Cells.Item(17, 8).Value = 9

Set-TextValue $q4 18 2 "009211"
$q4.Cells.Item(18, 3).Value = "中欧嘉和三年持有期混合C"
Set-TextValue $q4 18 4 "2.20"
Set-TextValue $q4 18 5 "93.03"
Set-TextValue $q4 18 6 "4.57"
Set-TextValue $q4 18 7 "0.1005"
$q4.Cells.Item(18, 8).Value = 2

Set-TextValue $q4 19 2 "013994"
$q4.Cells.Item(19, 3).Value = "中欧光熠一年持有期混合型证券投资基金C"
Set-TextValue $q4 19 4 "2.48"
Set-TextValue $q4 19 5 "87.20"
Set-TextValue $q4 19 6 "3.93"
Set-TextValue $q4 19 7 "0.0975"
$q4.Cells.Item(19, 8).Value = 3

Set-TextValue $q4 20 2 "011709"
$q4.Cells.Item(20, 3).Value = "中欧嘉益一年混合C"
Set-TextValue $q4 20 4 "1.82"
Set-TextValue $q4 20 5 "91.19"
Set-TextValue $q4 20 6 "4.23"
Set-TextValue $q4 20 7 "0.0770"
$q4.Cells.Item(20, 8).Value = 3

Set-TextValue $q4 21 2 "010724"
$q4.Cells.Item(21, 3).Value = "中欧价值成长混合C"
Set-TextValue $q4 21 4 "1.52"
Set-TextValue $q4 21 5 "83.01"
Set-TextValue $q4 21 6 "4.48"
Set-TextValue $q4 21 7 "0.0681"
$q4.Cells.Item(21, 8).Value = 4

Set-TextValue $q4 22 2 "000524"
$q4.Cells.Item(22, 3).Value = "上投摩根民生需求股票A"
Set-TextValue $q4 22 4 "1.32"
Set-TextValue $q4 22 5 "88.01"
Set-TextValue $q4 22 6 "3.13"
Set-TextValue $q4 22 7 "0.0413"
$q4.Cells.Item(22, 8).Value = 7

Set-TextValue $q4 23 2 "519678"
$q4.Cells.Item(23, 3).Value = "银河消费驱动混合A"
Set-TextValue $q4 23 4 "0.86"
Set-TextValue $q4 23 5 "83.68"
Set-TextValue $q4 23 6 "4.43"
Set-TextValue $q4 23 7 "0.0381"
$q4.Cells.Item(23, 8).Value = 10

Set-TextValue $q4 24 2 "010679"
$q4.Cells.Item(24, 3).Value = "中欧均衡成长混合C"
Set-TextValue $q4 24 4 "0.74"
Set-TextValue $q4 24 5 "83.93"
Set-TextValue $q4 24 6 "4.71"
Set-TextValue $q4 24 7 "0.0349"
$q4.Cells.Item(24, 8).Value = 4

Set-TextValue $q4 25 2 "012203"
$q4.Cells.Item(25, 3).Value = "中加消费优选混合C"
Set-TextValue $q4 25 4 "0.95"
Set-TextValue $q4 25 5 "88.36"
Set-TextValue $q4 25 6 "3.13"
Set-TextValue $q4 25 7 "0.0297"
$q4.Cells.Item(25, 8).Value = 9

Set-TextValue $q4 26 2 "010948"
$q4.Cells.Item(26, 3).Value = "中欧嘉选混合C"
Set-TextValue $q4 26 4 "0.74"
Set-TextValue $q4 26 5 "83.49"
Set-TextValue $q4 26 6 "3.22"
Set-TextValue $q4 26 7 "0.0238"
$q4.Cells.Item(26, 8).Value = 4

Set-TextValue $q4 27 2 "015032"
$q4.Cells.Item(27, 3).Value = "中融医药消费混合A"
Set-TextValue $q4 27 4 "0.52"
Set-TextValue $q4 27 5 "92.98"
Set-TextValue $q4 27 6 "3.53"
Set-TextValue $q4 27 7 "0.0184"
$q4.Cells.Item(27, 8).Value = 10

Set-TextValue $q4 28 2 "001883"
$q4.Cells.Item(28, 3).Value = "中欧新动力混合（LOF）E"
Set-TextValue $q4 28 4 "0.36"
Set-TextValue $q4 28 5 "86.77"
Set-TextValue $q4 28 6 "4.81"
Set-TextValue $q4 28 7 "0.0173"
$q4.Cells.Item(28, 8).Value = 3

Set-TextValue $q4 29 2 "000757"
$q4.Cells.Item(29, 3).Value = "华富智慧城市灵活配置混合"
Set-TextValue $q4 29 4 "0.48"
Set-TextValue $q4 29 5 "92.86"
Set-TextValue $q4 29 6 "3.24"
Set-TextValue $q4 29 7 "0.0156"
$q4.Cells.Item(29, 8).Value = 6

Set-TextValue $q4 30 2 "009565"
$q4.Cells.Item(30, 3).Value = "汇安消费龙头混合C"
Set-TextValue $q4 30 4 "0.45"
Set-TextValue $q4 30 5 "94.58"
Set-TextValue $q4 30 6 "3.27"
Set-TextValue $q4 30 7 "0.0147"
$q4.Cells.Item(30, 8).Value = 9

Set-TextValue $q4 31 2 "001482"
$q4.Cells.Item(31, 3).Value = "上投摩根新兴服务股票A"
Set-TextValue $q4 31 4 "0.35"
Set-TextValue $q4 31 5 "89.77"
Set-TextValue $q4 31 6 "3.09"
Set-TextValue $q4 31 7 "0.0108"
$q4.Cells.Item(31, 8).Value = 7

Set-TextValue $q4 32 2 "015086"
$q4.Cells.Item(32, 3).Value = "中欧核心消费股票C"
Set-TextValue $q4 32 4 "0.17"
Set-TextValue $q4 32 5 "84.11"
Set-TextValue $q4 32 6 "3.51"
Set-TextValue $q4 32 7 "0.0060"
$q4.Cells.Item(32, 8).Value = 10

Set-TextValue $q4 33 2 "015085"
$q4.Cells.Item(33, 3).Value = "中欧核心消费股票A"
Set-TextValue $q4 33 4 "0.15"
Set-TextValue $q4 33 5 "84.11"
Set-TextValue $q4 33 6 "3.51"
Set-TextValue $q4 33 7 "0.0053"
$q4.Cells.Item(33, 8).Value = 10

Set-TextValue $q4 34 2 "960016"
$q4.Cells.Item(34, 3).Value = "交银成长混合H"
Set-TextValue $q4 34 4 "0.17"
Set-TextValue $q4 34 5 "80.90"
Set-TextValue $q4 34 6 "2.56"
Set-TextValue $q4 34 7 "0.0044"
$q4.Cells.Item(34, 8).Value = 10

Set-TextValue $q4 35 2 "162211"
$q4.Cells.Item(35, 3).Value = "泰达宏利品质生活混合"
Set-TextValue $q4 35 4 "0.10"
Set-TextValue $q4 35 5 "76.66"
Set-TextValue $q4 35 6 "3.89"
Set-TextValue $q4 35 7 "0.0039"
$q4.Cells.Item(35, 8).Value = 6

Set-TextValue $q4 36 2 "007315"
$q4.Cells.Item(36, 3).Value = "汇安嘉盈一年持有期债券A"
Set-TextValue $q4 36 4 "0.16"
Set-TextValue $q4 36 5 "24.33"
Set-TextValue $q4 36 6 "1.24"
Set-TextValue $q4 36 7 "0.0020"
$q4.Cells.Item(36, 8).Value = 4

Set-TextValue $q4 37 2 "010270"
$q4.Cells.Item(37, 3).Value = "汇安嘉盈一年持有期债券C"
Set-TextValue $q4 37 4 "0.15"
Set-TextValue $q4 37 5 "24.33"
Set-TextValue $q4 37 6 "1.24"
Set-TextValue $q4 37 7 "0.0019"
$q4.Cells.Item(37, 8).Value = 4

Set-TextValue $q4 38 2 "015033"
$q4.Cells.Item(38, 3).Value = "中融医药消费混合C"
Set-TextValue $q4 38 4 "0.04"
Set-TextValue $q4 38 5 "92.98"
Set-TextValue $q4 38 6 "3.53"
Set-TextValue $q4 38 7 "0.0014"
$q4.Cells.Item(38, 8).Value = 10

Set-TextValue $q4 39 2 "015668"
$q4.Cells.Item(39, 3).Value = "银河消费驱动混合C"
Set-TextValue $q4 39 4 "0.01"
Set-TextValue $q4 39 5 "83.68"
Set-TextValue $q4 39 6 "4.43"
Set-TextValue $q4 39 7 "0.0004"
$q4.Cells.Item(39, 8).Value = 10

Set-TextValue $q4 40 2 "017099"
$q4.Cells.Item(40, 3).Value = "上投摩根民生需求股票C"
Set-TextValue $q4 40 4 "0.00"
Set-TextValue $q4 40 5 "88.01"
Set-TextValue $q4 40 6 "3.13"
$q4.Cells.Item(40, 7).Value = 0
$q4.Cells.Item(40, 8).Value = 7

Set-TextValue $q4 41 2 "017177"
$q4.Cells.Item(41, 3).Value = "上投摩根新兴服务股票C"
Set-TextValue $q4 41 4 "0.00"
Set-TextValue $q4 41 5 "89.77"
Set-TextValue $q4 41 6 "3.09"
$q4.Cells.Item(41, 7).Value = 0
$q4.Cells.Item(41, 8).Value = 7


# ---- 3. Insert the new 2022-Q4 row into "总计" --------------------------
# Shift the six existing data rows (currently rows 2-7) down to rows 3-8,
# copying values bottom-up so nothing is overwritten before it's read.
for ($r = 7; $r -ge 2; $r--) {
    $label = $zj.Cells.Item($r, 2).Value()
    $cnt   = $zj.Cells.Item($r, 3).Value()
    $mv    = $zj.Cells.Item($r, 4).Value()
    $zj.Cells.Item($r + 1, 1).Value = $r
    $zj.Cells.Item($r + 1, 2).Value = $label
    $zj.Cells.Item($r + 1, 3).Value = $cnt
    $zj.Cells.Item($r + 1, 4).Value = $mv
}

# Row 2 becomes the new 2022-Q4 entry.
$zj.Cells.Item(2, 1).Value = 0
$zj.Cells.Item(2, 2).Value = "2022-Q4"
$zj.Cells.Item(2, 3).Value = 40
$zj.Cells.Item(2, 4).Value = 7.52

# Give the new index cell (A8) the same bold+border style used by the rest
# of column A, by copying the format from the cell just above it.
$zj.Cells.Item(7, 1).Copy()
$zj.Cells.Item(8, 1).PasteSpecial(-4122)
